$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the "last updated" timestamp string (A1)
$ws.Range("A1").Value = "Datos actualizados a 16 de Agosto de 2020 a las 17:26"

# Update country data rows (values refreshed from the latest feed; three countries
# - Kenia/Moldavia, Trinidad yTobago/Tanzania, Montserrat/Islas Malvinas - swap rank order)

# Row 4: Estados Unidos
$ws.Range("A4").Value = "Estados Unidos"
$ws.Range("B4").Value = 5538962
$ws.Range("C4").Value = 9173
$ws.Range("D4").Value = 2904439
$ws.Range("E4").Value = 2461766
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 151
$ws.Range("H4").Value = 172757

# Row 6: India
$ws.Range("A6").Value = "India"
$ws.Range("B6").Value = 2634256
$ws.Range("C6").Value = 45048
$ws.Range("D6").Value = 1904612
$ws.Range("E6").Value = 678799
$ws.Range("F6").Value = 0
$ws.Range("G6").Value = 761
$ws.Range("H6").Value = 50845

# Row 15: Reino Unido
$ws.Range("A15").Value = "Reino Unido"
$ws.Range("B15").Value = 318484
$ws.Range("C15").Value = 1040
$ws.Range("D15").Value = 0
$ws.Range("E15").Value = 0
$ws.Range("F15").Value = 0
$ws.Range("G15").Value = 5
$ws.Range("H15").Value = 41366

# Row 38: Oman
$ws.Range("A38").Value = "Oman"
$ws.Range("B38").Value = 83086
$ws.Range("C38").Value = 162
$ws.Range("D38").Value = 77680
$ws.Range("E38").Value = 4834
$ws.Range("F38").Value = 0
$ws.Range("G38").Value = 10
$ws.Range("H38").Value = 572

# Row 48: Singapur
$ws.Range("A48").Value = "Singapur"
$ws.Range("B48").Value = 55747
$ws.Range("C48").Value = 86
$ws.Range("D48").Value = 51953
$ws.Range("E48").Value = 3767
$ws.Range("F48").Value = 0
$ws.Range("G48").Value = 0
$ws.Range("H48").Value = 27

# Row 64: Kenia
$ws.Range("A64").Value = "Kenia"
$ws.Range("B64").Value = 30120
$ws.Range("C64").Value = 271
$ws.Range("D64").Value = 16656
$ws.Range("E64").Value = 12990
$ws.Range("F64").Value = 0
$ws.Range("G64").Value = 2
$ws.Range("H64").Value = 474

# Row 65: Moldavia
$ws.Range("A65").Value = "Moldavia"
$ws.Range("B65").Value = 29905
$ws.Range("C65").Value = 0
$ws.Range("D65").Value = 20908
$ws.Range("E65").Value = 8102
$ws.Range("F65").Value = 0
$ws.Range("G65").Value = 0
$ws.Range("H65").Value = 895

# Row 95: Tayikistan
$ws.Range("A95").Value = "Tayikistan"
$ws.Range("B95").Value = 8065
$ws.Range("C95").Value = 36
$ws.Range("D95").Value = 6855
$ws.Range("E95").Value = 1146
$ws.Range("F95").Value = 0
$ws.Range("G95").Value = 0
$ws.Range("H95").Value = 64

# Row 119: Somalia
$ws.Range("A119").Value = "Somalia"
$ws.Range("B119").Value = 3256
$ws.Range("C119").Value = 6
$ws.Range("D119").Value = 2374
$ws.Range("E119").Value = 789
$ws.Range("F119").Value = 0
$ws.Range("G119").Value = 0
$ws.Range("H119").Value = 93

# Row 166: Trinidad yTobago
$ws.Range("A166").Value = "Trinidad yTobago"
$ws.Range("B166").Value = 519
$ws.Range("C166").Value = 22
$ws.Range("D166").Value = 140
$ws.Range("E166").Value = 368
$ws.Range("F166").Value = 0
$ws.Range("G166").Value = 1
$ws.Range("H166").Value = 11

# Row 167: Tanzania
$ws.Range("A167").Value = "Tanzania"
$ws.Range("B167").Value = 509
$ws.Range("C167").Value = 0
$ws.Range("D167").Value = 183
$ws.Range("E167").Value = 305
$ws.Range("F167").Value = 0
$ws.Range("G167").Value = 0
$ws.Range("H167").Value = 21

# Row 213: Montserrat
$ws.Range("A213").Value = "Montserrat"
$ws.Range("B213").Value = 13
$ws.Range("C213").Value = 0
$ws.Range("D213").Value = 12
$ws.Range("E213").Value = 0
$ws.Range("F213").Value = 0
$ws.Range("G213").Value = 0
$ws.Range("H213").Value = 1

# Row 214: Islas Malvinas
$ws.Range("A214").Value = "Islas Malvinas"
$ws.Range("B214").Value = 13
$ws.Range("C214").Value = 0
$ws.Range("D214").Value = 13
$ws.Range("E214").Value = 0
$ws.Range("F214").Value = 0
$ws.Range("G214").Value = 0
$ws.Range("H214").Value = 0
